## Edit: (1) change the table style on the table in slide 16, and
## (2) swap the presentation's theme colour palette from the custom
## "Integral" scheme to the standard "Office" scheme (dk1/lt1 stay the
## same black/white; the other 10 theme colours change).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Table style change (slide 16 -> the graphicFrame holding a table)
# ---------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{89843619-0EDF-4FBD-9E9A-B3DB7476ECB2}")
    }
}

# ---------------------------------------------------------------
# 2) Theme colour scheme change (Integral -> Office)
#    ThemeColorScheme index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
#    5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5,
#    10 accent6, 11 hlink, 12 folHlink
# ---------------------------------------------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme

# RGB() long value = R + G*256 + B*65536 (standard VBA colour encoding)
$tcs.Item(3).RGB  = 6968388    # dk2      455F51 -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E3DED1 -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  99CB38 -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  63A537 -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  E6D024 -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  CC9700 -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4EB3CF -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  378DA6 -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    6B9F25 -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink B26B02 -> 954F72
